# patch 1.12 | oprava responzivity pri navigacii
#
# Diel sheet: the three-row "Kvant / Skenerovy drziak KV30-cast X /
# 99-23-3426" block (rows 2-4) and the two-row "Slavia / FRAESTEIL /
# 3B8B034_0209" block (rows 5-6) had the firma/nazov/cisloVykresu
# (A:C) and upnutie (H, only where it repeats) values copy-pasted on
# every row. Turn those repeats into merged cells (center/middle
# aligned) instead, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diel")

$xlCenter = -4108

function Merge-Centered($ws, $rangeAddress) {
    # Align first, then merge *that exact range object* - merging a
    # wider/previously-touched range here makes later merges in the
    # same run pick up a stale effective style.
    $r = $ws.Range($rangeAddress)
    $r.HorizontalAlignment = $xlCenter
    $r.VerticalAlignment = $xlCenter
    $r.Merge()
}

# Group 1 (rows 2-4): firma / nazov / cisloVykresu repeated 3x, upnutie
# repeated only on rows 2-3 (row 4 has a different upnutie value).
Merge-Centered $ws "A2:A4"
Merge-Centered $ws "B2:B4"
Merge-Centered $ws "C2:C4"
Merge-Centered $ws "H2:H3"

# Group 2 (rows 5-6): firma / nazov / cisloVykresu + upnutie repeated 2x.
Merge-Centered $ws "A5:A6"
Merge-Centered $ws "B5:B6"
Merge-Centered $ws "C5:C6"
Merge-Centered $ws "H5:H6"

# Move the active selection (navigation responsiveness fix).
$ws.Range("J7").Select()

Write-Host "Applied merges + selection change to Diel sheet"
